$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 21
$ws.Range("H21").Value = 21533.334
$ws.Range("J21").Value = 21533.334
$ws.Range("L21").Value = 21533.334
$ws.Range("N21").Value = -22469.334
# Row 23
$ws.Range("H23").Value = 21533.334
$ws.Range("J23").Value = 21533.334
$ws.Range("L23").Value = 21533.334
$ws.Range("N23").Value = -22001.334
# Row 33
$ws.Range("H33").Value = 196.78947
$ws.Range("I33").Value = 222.4375
$ws.Range("K33").Value = 222.4375
$ws.Range("M33").Value = 6.5625
# Row 51
$ws.Range("H51").Value = 6200
$ws.Range("I51").Value = 5666.6665
$ws.Range("K51").Value = 5666.6665
$ws.Range("L51").Value = 7000
$ws.Range("M51").Value = -5182.6665
$ws.Range("N51").Value = -7968
# Row 53
$ws.Range("H53").Value = 305.5
$ws.Range("I53").Value = 341.7
$ws.Range("J53").Value = 215
$ws.Range("K53").Value = 341.7
$ws.Range("L53").Value = 215
$ws.Range("M53").Value = 295.3
$ws.Range("N53").Value = -1489
# Row 86
$ws.Range("H86").Value = 5594.174
$ws.Range("I86").Value = 1251.2307
$ws.Range("J86").Value = 11240
$ws.Range("K86").Value = 1251.2307
$ws.Range("L86").Value = 11240
$ws.Range("M86").Value = -128.2307000000001
$ws.Range("N86").Value = -13486
# Row 89
$ws.Range("H89").Value = 5594.174
$ws.Range("I89").Value = 1251.2307
$ws.Range("J89").Value = 11240
$ws.Range("K89").Value = 6256.1535
$ws.Range("L89").Value = 56200
$ws.Range("M89").Value = -640.1535000000003
$ws.Range("N89").Value = -67432
# Row 100
$ws.Range("H100").Value = 1353.6
$ws.Range("I100").Value = 800.625
$ws.Range("J100").Value = 1985.5714
$ws.Range("K100").Value = 800.625
$ws.Range("L100").Value = 1985.5714
$ws.Range("M100").Value = -259.625
$ws.Range("N100").Value = -3067.5714
# Row 129
$ws.Range("H129").Value = 754.63635
$ws.Range("J129").Value = 800.42
$ws.Range("L129").Value = 2401.26
$ws.Range("N129").Value = -12401.26

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6240.3164
$ws.Range("I32").Value = 4815.1357
$ws.Range("J32").Value = 13030.883
$ws.Range("K32").Value = 4815.1357
$ws.Range("L32").Value = 13030.883
$ws.Range("M32").Value = -4528.1357
$ws.Range("N32").Value = -13604.883
# Row 74
$ws.Range("H74").Value = 30304558
$ws.Range("I74").Value = 40000656
$ws.Range("J74").Value = 4247.5
$ws.Range("K74").Value = 40000656
$ws.Range("L74").Value = 4247.5
$ws.Range("M74").Value = -39999782
$ws.Range("N74").Value = -5995.5
# Row 77
$ws.Range("H77").Value = 30304558
$ws.Range("I77").Value = 40000656
$ws.Range("J77").Value = 4247.5
$ws.Range("K77").Value = 200003280
$ws.Range("L77").Value = 21237.5
$ws.Range("M77").Value = -199998912
$ws.Range("N77").Value = -29973.5
# Row 132
$ws.Range("H132").Value = 15212.359
$ws.Range("I132").Value = 2201.742
$ws.Range("J132").Value = 65628.5
$ws.Range("K132").Value = 6605.226000000001
$ws.Range("L132").Value = 196885.5
$ws.Range("M132").Value = -4075.226000000001
$ws.Range("N132").Value = -201945.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 744.6087
$ws.Range("I94").Value = 601.44446
$ws.Range("J94").Value = 1260
$ws.Range("K94").Value = 601.44446
$ws.Range("L94").Value = 1260
$ws.Range("M94").Value = -150.44446
$ws.Range("N94").Value = -2162

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3281.7302
$ws.Range("I31").Value = 1509.7273
$ws.Range("J31").Value = 5230.933
$ws.Range("K31").Value = 1509.7273
$ws.Range("L31").Value = 5230.933
$ws.Range("M31").Value = -1214.7273
$ws.Range("N31").Value = -5820.933
# Row 34
$ws.Range("H34").Value = 3281.7302
$ws.Range("I34").Value = 1509.7273
$ws.Range("J34").Value = 5230.933
$ws.Range("K34").Value = 1509.7273
$ws.Range("L34").Value = 5230.933
$ws.Range("M34").Value = -1307.7273
$ws.Range("N34").Value = -5634.933
# Row 105
$ws.Range("H105").Value = 745.4167
$ws.Range("I105").Value = 745.4167
$ws.Range("K105").Value = 745.4167
$ws.Range("M105").Value = 1001.5833

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 82
$ws.Range("H82").Value = 9999.166999999999
$ws.Range("J82").Value = 9999.166999999999
$ws.Range("L82").Value = 29997.501
$ws.Range("N82").Value = -30809.501
# Row 85
$ws.Range("H85").Value = 9999.166999999999
$ws.Range("J85").Value = 9999.166999999999
$ws.Range("L85").Value = 29997.501
$ws.Range("N85").Value = -32805.501
# Row 113
$ws.Range("H113").Value = 705.3182
$ws.Range("J113").Value = 706.5
$ws.Range("L113").Value = 2119.5
$ws.Range("N113").Value = -6459.5
# Row 131
$ws.Range("H131").Value = 736.74
$ws.Range("J131").Value = 749.96906
$ws.Range("L131").Value = 2249.90718
$ws.Range("N131").Value = -12329.90718

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3317.5833
$ws.Range("I80").Value = 3052
$ws.Range("J80").Value = 3542.3076
$ws.Range("K80").Value = 3052
$ws.Range("L80").Value = 3542.3076
$ws.Range("M80").Value = -2054
$ws.Range("N80").Value = -5538.3076
# Row 83
$ws.Range("H83").Value = 3317.5833
$ws.Range("I83").Value = 3052
$ws.Range("J83").Value = 3542.3076
$ws.Range("K83").Value = 15260
$ws.Range("L83").Value = 17711.538
$ws.Range("M83").Value = -10268
$ws.Range("N83").Value = -27695.538
# Row 126
$ws.Range("H126").Value = 3016.6667
$ws.Range("I126").Value = 2127.6956
$ws.Range("J126").Value = 3746.8928
$ws.Range("K126").Value = 6383.0868
$ws.Range("L126").Value = 11240.6784
$ws.Range("M126").Value = -3913.0868
$ws.Range("N126").Value = -16180.6784

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4435.037
$ws.Range("I7").Value = 4312.3
$ws.Range("J7").Value = 4785.7144
$ws.Range("K7").Value = 4312.3
$ws.Range("L7").Value = 4785.7144
$ws.Range("M7").Value = -4200.3
$ws.Range("N7").Value = -5009.7144
# Row 22
$ws.Range("H22").Value = 4669.6113
$ws.Range("I22").Value = 5070.2
$ws.Range("J22").Value = 2666.6667
$ws.Range("K22").Value = 5070.2
$ws.Range("L22").Value = 2666.6667
$ws.Range("M22").Value = -4775.2
$ws.Range("N22").Value = -3256.6667
# Row 27
$ws.Range("H27").Value = 4669.6113
$ws.Range("I27").Value = 5070.2
$ws.Range("J27").Value = 2666.6667
$ws.Range("K27").Value = 5070.2
$ws.Range("L27").Value = 2666.6667
$ws.Range("M27").Value = -4963.2
$ws.Range("N27").Value = -2880.6667
# Row 93
$ws.Range("H93").Value = 1508.84
$ws.Range("I93").Value = 1391.5238
$ws.Range("J93").Value = 2124.75
$ws.Range("K93").Value = 1391.5238
$ws.Range("L93").Value = 2124.75
$ws.Range("M93").Value = -143.5237999999999
$ws.Range("N93").Value = -4620.75
# Row 100
$ws.Range("H100").Value = 1960.8235
$ws.Range("I100").Value = 1052.5714
$ws.Range("J100").Value = 2596.6
$ws.Range("K100").Value = 1052.5714
$ws.Range("L100").Value = 2596.6
$ws.Range("M100").Value = -511.5714
$ws.Range("N100").Value = -3678.6
# Row 126
$ws.Range("H126").Value = 4435.037
$ws.Range("I126").Value = 4312.3
$ws.Range("J126").Value = 4785.7144
$ws.Range("K126").Value = 12936.9
$ws.Range("L126").Value = 14357.1432
$ws.Range("M126").Value = -10466.9
$ws.Range("N126").Value = -19297.1432

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 4537
$ws.Range("I62").Value = 2611
$ws.Range("J62").Value = 5500
$ws.Range("K62").Value = 2611
$ws.Range("L62").Value = 5500
$ws.Range("M62").Value = -1987
$ws.Range("N62").Value = -6748
# Row 65
$ws.Range("H65").Value = 4537
$ws.Range("I65").Value = 2611
$ws.Range("J65").Value = 5500
$ws.Range("K65").Value = 13055
$ws.Range("L65").Value = 27500
$ws.Range("M65").Value = -9935
$ws.Range("N65").Value = -33740
# Row 107
$ws.Range("H107").Value = 142857440
$ws.Range("I107").Value = 166666960
$ws.Range("J107").Value = 290
$ws.Range("K107").Value = 500000880
$ws.Range("L107").Value = 870
$ws.Range("M107").Value = -499998960
$ws.Range("N107").Value = -4710
# Row 136
$ws.Range("H136").Value = 30363838
$ws.Range("I136").Value = 41291864
$ws.Range("K136").Value = 123875592
$ws.Range("M136").Value = -123873042
